$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 6 (pushes AttributeType enum items down by one),
# to hold the new BEHAVIOR_TYPE enum item.
$ws.Rows("6:6").Insert()

# New row 6: BEHAVIOR_TYPE = 2
$ws.Range("G6").Value = "BEHAVIOR_TYPE"
$ws.Range("I6").Value = 2

# Renumber the AttributeType values (ATTACK..PIERCE), now at rows 7-14,
# from 2-9 to 10-17 since BEHAVIOR_TYPE took value 2.
$ws.Range("I7").Value = 10
$ws.Range("I8").Value = 11
$ws.Range("I9").Value = 12
$ws.Range("I10").Value = 13
$ws.Range("I11").Value = 14
$ws.Range("I12").Value = 15
$ws.Range("I13").Value = 16
$ws.Range("I14").Value = 17

# Match the author's final selection.
$ws.Range("I14").Select() | Out-Null
